$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E corresponds to the "municipio-nombre" dimension metadata block
# (rows 2-4 hold iaest/sdmx dimension metadata: type, class, URI template).
# Re-process with the newly curated dimensions: municipio-nombre now
# mirrors the refArea / dim / URI-Municipio pattern used by the other
# area columns (provincia-nombre, comarca-nombre).
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
